# Apply the "updated scrape" cell edits to the 展览 (sheet1) and 全部类型 (sheet4)
# worksheets. Both sheets carry the same underlying rows (全部类型 has a few
# extra rows mixed in from 演出), so we update each sheet independently using
# its own row numbers.

$wb = $excel.ActiveWorkbook

# Edits for sheet "展览" (rows 2-25)
$sheetExpo = $wb.Worksheets.Item("展览")
$expoEdits = @(
    @{ Ref = "G2";  Value = 60 },
    @{ Ref = "F3";  Value = 1816 },
    @{ Ref = "G3";  Value = 60 },
    @{ Ref = "F4";  Value = 41 },
    @{ Ref = "G5";  Value = "不可售" },
    @{ Ref = "F6";  Value = 676 },
    @{ Ref = "G6";  Value = 55 },
    @{ Ref = "G7";  Value = 60 },
    @{ Ref = "G8";  Value = 50 },
    @{ Ref = "G9";  Value = 55 },
    @{ Ref = "F13"; Value = 172 },
    @{ Ref = "F14"; Value = 27 },
    @{ Ref = "F17"; Value = 112 },
    @{ Ref = "F18"; Value = 5190 },
    @{ Ref = "F22"; Value = 2303 },
    @{ Ref = "F23"; Value = 75 },
    @{ Ref = "F24"; Value = 33 },
    @{ Ref = "F25"; Value = 2145 }
)
foreach ($edit in $expoEdits) {
    $sheetExpo.Range($edit.Ref).Value = $edit.Value
}

# Edits for sheet "全部类型" (rows 2-28; row numbers differ from 展览 because
# this sheet interleaves rows from 演出 as well)
$sheetAll = $wb.Worksheets.Item("全部类型")
$allEdits = @(
    @{ Ref = "G2";  Value = 60 },
    @{ Ref = "F3";  Value = 1816 },
    @{ Ref = "G3";  Value = 60 },
    @{ Ref = "F4";  Value = 41 },
    @{ Ref = "G5";  Value = "不可售" },
    @{ Ref = "F6";  Value = 676 },
    @{ Ref = "G6";  Value = 55 },
    @{ Ref = "G7";  Value = 60 },
    @{ Ref = "G8";  Value = 50 },
    @{ Ref = "G9";  Value = 55 },
    @{ Ref = "F13"; Value = 172 },
    @{ Ref = "F14"; Value = 27 },
    @{ Ref = "F17"; Value = 112 },
    @{ Ref = "F18"; Value = 5190 },
    @{ Ref = "F24"; Value = 2303 },
    @{ Ref = "F25"; Value = 75 },
    @{ Ref = "F27"; Value = 33 },
    @{ Ref = "F28"; Value = 2145 }
)
foreach ($edit in $allEdits) {
    $sheetAll.Range($edit.Ref).Value = $edit.Value
}
